$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Lay down row 3 using row 2's current (pre-edit) formatting, so
#    every cell in the new row inherits the right border/fill style
#    before anything gets its one-off "unlocked" look.
# ------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Helper-ish steps: write text that *looks* numeric ("10031", "10024")
# as real text (shared string), without leaving stray number-format
# styles behind. Recipe: mark the cell as Text, type the value, strip
# the format back off, then restore the plain bordered look by
# pasting formats from a same-style neighbour (C2, which is style 3
# and never changes).
# ------------------------------------------------------------------

# --- Row 2: Budget ID 10024 -> 10031, Check string updated ---
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "10031"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 5).Value = "UBS/JBBGT/10031/11/2023"

# Give F2 its one-off "unlocked / white fill" look.
$ws.Cells.Item(2, 6).Interior.Color = 16777215
$ws.Cells.Item(2, 6).Locked = $false

# --- Row 3: new JBBMT row ---
$ws.Cells.Item(3, 1).Value = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(3, 2).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(3, 3).Value = "JBBMT"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "10024"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(3, 5).Value = "UBS/JBBMT/10024/11/2023"
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Formula = "=SUM(F3:H3)"

# Give G3 the same one-off "unlocked / white fill" look that F2 has.
$ws.Cells.Item(3, 7).Interior.Color = 16777215
$ws.Cells.Item(3, 7).Locked = $false

# ------------------------------------------------------------------
# 2) Column E needs to be a touch wider.
# ------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 24.7109375

# ------------------------------------------------------------------
# 3) Mirror the conditional-format rule from I2 onto the new I3 cell.
# ------------------------------------------------------------------
$rule = $ws.Range("I3").FormatConditions.Add(1, 5, "=100")
$rule.Interior.Color = 255

# ------------------------------------------------------------------
# 4) Protect the worksheet (matches <sheetProtection sheet objects
#    scenarios/>); per-cell Locked state set above decides which
#    cells stay editable.
# ------------------------------------------------------------------
$ws.Protect("", $true, $true, $true)
